# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the newly scraped data, as described by the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F.
$updates = @{
    2  = 311
    3  = 102
    4  = 395
    5  = 11637
    6  = 824
    10 = 147
    11 = 169
    12 = 23
    13 = 49
    17 = 334
    18 = 1358
    19 = 77
    20 = 903
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
